$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# "View All" — the placeholder row that used to sit at row 49 (a one-off
# row formatted with the extra Arial / Malgun Gothic / Inconsolata fonts)
# is removed, so the everyday-language rows beneath it shift up by one.
$ws.Rows.Item(49).Delete()

# Reflect the viewport the author ended on after the edit.
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Range("B79").Select()
